$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "analysis_B4_BT02_09.56_10.40.xlsx"
$ws.Range("C1").Value = "analysis_B4_BT02_13.59_14.50.xlsx"
$ws.Range("D1").Value = "analysis_B4_BT03_14.54_15.49.xlsx"
$ws.Range("E1").Value = "analysis_B4_BT04_10.45_11.45.xlsx"
$ws.Range("F1").Value = "analysis_B4_BT04_15.52_16.57.xlsx"
$ws.Range("G1").Value = "analysis_B4_BT05_11.47_12.42.xlsx"
$ws.Range("H1").Value = "analysis_B4_BT06_09.00_09.55.xlsx"
$ws.Range("I1").Value = "analysis_B4_BT06_13.03_13.56.xlsx"
$ws.Range("J1").Value = "analysis_B4_BT06_17.05_17.45.xlsx"
